$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the boolean value in D3 (was TRUE, now empty)
$ws.Range("D3").ClearContents()

# Update the active selection to F10 (side-effect of user's navigation)
$ws.Range("F10").Select()
